# 014 Week 8 data update
# Fills in "WK 8" results for Sheet1 (Sunday pairs) column K and
# THURSDAY SINGLES column I. Dependent SUM()/COUNTIF() formulas on
# Sheet1, xxDO NOT EDITxx, etc. recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- Sheet1 : WK 8 column (K) -------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$sheet1Updates = @{
    12 = 32.0
    15 = 36.0
    16 = 23.0
    17 = 31.0
    20 = 37.0
    21 = 42.0
    22 = 35.0
    25 = 31.0
    26 = 30.0
    27 = 39.0
    28 = 32.0
    29 = 38.0
    31 = 31.0
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("K$row").Value = $sheet1Updates[$row]
}

# --- THURSDAY SINGLES : WK 8 column (I) ----------------------------------
$ws2 = $wb.Worksheets.Item("THURSDAY SINGLES")

$sheet2Updates = @{
    5  = 29.0
    7  = 29.0
    8  = 33.0
    12 = 34.0
    16 = 27.0
}

foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("I$row").Value = $sheet2Updates[$row]
}

$excel.CalculateFull()
